$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2023-01-09 12:56:41"
$newTimestamp = "2023-01-09 14:17:57"

for ($r = 2; $r -le 398; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

$mCell = $ws.Cells.Item(281, 13)
$mCell.Value = "Betty Bossi Pfaffenhut 2x  100g - Online kein Bestand 4.60 Schweizer Franken"
